$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.015.70'
$ws.Range("E2").Value = '  -1.35%  '
$ws.Range("D3").Value = '3.409.62'
$ws.Range("E3").Value = '  -1.79%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'576.53"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").Value = "'148.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = "'0.484"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.82%  '
$ws.Range("D9").Value = "'8.00"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.99%  '
$ws.Range("E10").Value = '  -0.86%  '
$ws.Range("E11").Value = '  +3.21%  '
$ws.Range("D12").Value = '3.993.89'
$ws.Range("E12").Value = '  -1.68%  '
$ws.Range("E13").Value = '  -0.44%  '
$ws.Range("D14").Value = "'28.28"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -5.13%  '
$ws.Range("D15").Value = '3.410.20'
$ws.Range("E15").Value = '  -1.52%  '
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("D17").Value = '61.938.80'
$ws.Range("E17").Value = '  -1.45%  '
$ws.Range("D18").Value = "'6.39"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.39%  '
$ws.Range("D19").Value = "'14.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.77%  '
$ws.Range("D20").Value = "'8.93"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.94%  '
$ws.Range("D21").Value = "'380.17"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -2.35%  '
$ws.Range("D22").Value = "'0.566"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +1.70%  '
$ws.Range("D23").Value = "'74.73"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("E24").Value = '  -0.18%  '
$ws.Range("D25").Value = '3.575.60'
$ws.Range("E25").Value = '  -0.77%  '
$ws.Range("E26").Value = '  -3.01%  '
$ws.Range("E27").Value = '  +1.12%  '
$ws.Range("D28").Value = "'7.65"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = "'7.92"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.59%  '
$ws.Range("E31").Value = '  -0.19%  '
$ws.Range("E32").Value = '  +0.00%  '
$ws.Range("E33").Value = '  -2.84%  '
$ws.Range("D34").Value = "'23.11"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.13%  '
$ws.Range("D35").Value = "'5.43"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.51%  '
$ws.Range("D36").Value = "'1.63"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +3.44%  '
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = "'6.90"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.12%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").Value = "'169.06"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.01%  '
$ws.Range("D39").Value = "'30.63"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.56%  '
$ws.Range("D40").Value = '3.442.42'
$ws.Range("E40").Value = '  -1.83%  '
$ws.Range("E41").Value = '  +3.45%  '
$ws.Range("D42").Value = "'0.783"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.03%  '
$ws.Range("D43").Value = "'42.33"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("D44").Value = "'4.37"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.87%  '
$ws.Range("E45").Value = '  -2.16%  '
$ws.Range("E46").Value = '  -3.39%  '
$ws.Range("D47").Value = '2.541.41'
$ws.Range("E47").Value = '  -2.46%  '
$ws.Range("E48").Value = '  +2.90%  '
$ws.Range("E49").Value = '  -3.65%  '
$ws.Range("D50").Value = "'22.55"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.36%  '
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +0.02%  '

Write-Host "Applied all changes"